$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 551.0680923226339
$ws.Range("D2").Value = 135.98277249591
$ws.Range("F2").Value = 448
$ws.Range("G2").Value = 507
$ws.Range("H2").Value = 619
$ws.Range("C3").Value = 40.30310759906428
$ws.Range("D3").Value = 5.324267421150879
$ws.Range("F3").Value = 37.17
$ws.Range("G3").Value = 40.02
$ws.Range("H3").Value = 43.41
$ws.Range("C4").Value = 1.620580493539557
$ws.Range("D4").Value = 2.151637427293199
$ws.Range("G4").Value = 1.08
$ws.Range("H4").Value = 2.05
$ws.Range("C5").Value = 323.2946474308017
$ws.Range("D5").Value = 11.14856390329808
$ws.Range("F5").Value = 315.94
$ws.Range("G5").Value = 324.65
$ws.Range("H5").Value = 332.2
$ws.Range("C6").Value = 20.78274570178028
$ws.Range("D6").Value = 2.53157948035533
$ws.Range("F6").Value = 19.41
$ws.Range("G6").Value = 20.79
$ws.Range("H6").Value = 22.25
$ws.Range("C7").Value = -76.06973622877116
$ws.Range("D7").Value = 22.89142367125974
$ws.Range("G7").Value = -72
$ws.Range("C8").Value = 7.644266129774256
$ws.Range("D8").Value = 6.897221567992845
$ws.Range("C9").Value = 9.321485312455726
$ws.Range("D9").Value = 1.685235526711948
$ws.Range("C10").Value = 867.8301709770533
$ws.Range("D10").Value = 0.4614231124990945
$ws.Range("C11").Value = 0.5554448833875715
$ws.Range("D11").Value = 0.5887632158809544
$ws.Range("C12").Value = 22.73994378273617
$ws.Range("D12").Value = 12.29186665115116
$ws.Range("C13").Value = 0.6738674081548559
$ws.Range("D13").Value = 0.7505400353874433
$ws.Range("C14").Value = 1.826958330223971
$ws.Range("D14").Value = 1.664121300141535
$ws.Range("C15").Value = 93.46973622877101
$ws.Range("D15").Value = 22.89142367125974
$ws.Range("G15").Value = 89.40000000000001
$ws.Range("C16").Value = -85.31226661389034
$ws.Range("D16").Value = 20.60246642824096
$ws.Range("F16").Value = -101.8707776445072
$ws.Range("H16").Value = -67.79706163635328
$ws.Range("C17").Value = -77.66800048411609
$ws.Range("D17").Value = 25.30226435379067
$ws.Range("F17").Value = -92.66683163887967
$ws.Range("G17").Value = -72.71081852649533
$ws.Range("H17").Value = -57.22214159641585
